# Auto-generated edit script: refresh crypto price/volume data
# (and two pairs of row-swaps) to match the Fri Jun 23 13:54:42 UTC 2023
# GitHub Actions data refresh commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '30.041.35'
$ws.Range("E2").Value = '  -0.64%  '
# Row 3
$ws.Range("D3").Value = '1.870.66'
$ws.Range("E3").Value = '  -1.37%  '
# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.002'
$ws.Range("E4").Value = '  +0.57%  '
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '241.75'
$ws.Range("E5").Value = '  -2.39%  '
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.001'
$ws.Range("E6").Value = '  +0.41%  '
# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4866'
$ws.Range("E7").Value = '  -2.54%  '
# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2887'
$ws.Range("E8").Value = '  -2.10%  '
# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06556'
$ws.Range("E9").Value = '  -1.70%  '
# Row 10
$ws.Range("D10").Value = '1.875.45'
$ws.Range("E10").Value = '  -1.10%  '
# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '16.29'
$ws.Range("E11").Value = '  -4.36%  '
# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07180'
$ws.Range("E12").Value = '  -0.33%  '
# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.6606'
$ws.Range("E13").Value = '  -2.69%  '
# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.884'
$ws.Range("E14").Value = '  +0.45%  '
# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '85.47'
$ws.Range("E15").Value = '  -0.65%  '
# Row 16
$ws.Range("D16").Value = '30.053.45'
$ws.Range("E16").Value = '  -0.45%  '
# Row 17
$ws.Range("B17").Value = 'Dai'
$ws.Range("C17").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.001'
$ws.Range("E17").Value = '  +0.24%  '
# Row 18
$ws.Range("B18").Value = 'ShibaInu'
$ws.Range("C18").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000007741'
$ws.Range("E18").Value = '  -3.47%  '
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.68'
$ws.Range("E19").Value = '  -1.94%  '
# Row 20
$ws.Range("D20").Value = '2.118.72'
$ws.Range("E20").Value = '  -0.75%  '
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.002'
$ws.Range("E21").Value = '  +0.59%  '
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.729'
$ws.Range("E22").Value = '  -1.07%  '
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.798'
$ws.Range("E23").Value = '  +2.49%  '
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.145'
$ws.Range("E24").Value = '  -0.29%  '
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '152.00'
$ws.Range("E25").Value = '  +3.21%  '
# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '142.46'
$ws.Range("E26").Value = '  +6.46%  '
# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '16.88'
$ws.Range("E27").Value = '  +0.13%  '
# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.864'
$ws.Range("E28").Value = '  -4.53%  '
# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.396'
$ws.Range("E29").Value = '  +1.38%  '
# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.178'
$ws.Range("E30").Value = '  -1.30%  '
# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08731'
$ws.Range("E31").Value = '  -0.26%  '
# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.968'
$ws.Range("E32").Value = '  +0.19%  '
# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05100'
$ws.Range("E33").Value = '  -1.19%  '
# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7050'
$ws.Range("E34").Value = '  -0.74%  '
# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.099'
$ws.Range("E35").Value = '  -2.15%  '
# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.667'
$ws.Range("E36").Value = '  +0.11%  '
# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.01821'
$ws.Range("E37").Value = '  +9.31%  '
# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.672'
$ws.Range("E38").Value = '  -4.03%  '
# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.111'
$ws.Range("E39").Value = '  -5.93%  '
# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.9217'
$ws.Range("E40").Value = '  -2.09%  '
# Row 41
$ws.Range("B41").Value = 'PaxDollar'
$ws.Range("C41").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.9990'
$ws.Range("E41").Value = '  +0.21%  '
# Row 42
$ws.Range("B42").Value = 'Quant'
$ws.Range("C42").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '103.38'
$ws.Range("E42").Value = '  +0.15%  '
# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.731'
$ws.Range("E43").Value = '  -6.00%  '
# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.4187'
$ws.Range("E44").Value = '  -0.95%  '
# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '7.364'
$ws.Range("E45").Value = '  -2.09%  '
# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.1273'
$ws.Range("E46").Value = '  +0.22%  '
# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.05708'
$ws.Range("E47").Value = '  -0.29%  '
# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '32.66'
$ws.Range("E48").Value = '  -0.66%  '
# Row 49
$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.248'
$ws.Range("E49").Value = '  +0.09%  '
# Row 50
$ws.Range("B50").Value = 'Decentraland'
$ws.Range("C50").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.3726'
$ws.Range("E50").Value = '  -0.62%  '
# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '55.49'
$ws.Range("E51").Value = '  -1.12%  '
